$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "68.825.98"
$ws.Range("E2").Value = "  +1.93%  "

$ws.Range("D3").Value = "3.281.59"
$ws.Range("E3").Value = "  +0.88%  "

$ws.Range("E4").Value = "  -0.01%  "

$ws.Range("D5").Value = "'585.17"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.12%  "

$ws.Range("D6").Value = "'182.56"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.56%  "

$ws.Range("E7").Value = "  -0.03%  "

$ws.Range("E8").Value = "  +0.98%  "

$ws.Range("E9").Value = "  +3.38%  "

$ws.Range("D10").Value = "'6.66"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.70%  "

$ws.Range("E11").Value = "  +1.89%  "

$ws.Range("D12").Value = "3.858.66"
$ws.Range("E12").Value = "  +0.94%  "

$ws.Range("E13").Value = "  -0.23%  "

$ws.Range("D14").Value = "'28.85"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.24%  "

$ws.Range("D15").Value = "68.783.32"
$ws.Range("E15").Value = "  +1.91%  "

$ws.Range("D16").Value = "'0.0000172"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +2.87%  "

$ws.Range("D17").Value = "3.270.88"
$ws.Range("E17").Value = "  +0.51%  "

$ws.Range("D18").Value = "'5.83"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.35%  "

$ws.Range("D19").Value = "'13.60"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.61%  "

$ws.Range("D20").Value = "'394.65"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +4.66%  "

$ws.Range("D21").Value = "'7.73"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.56%  "

$ws.Range("D22").Value = "'71.56"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.18%  "

$ws.Range("D23").Value = "'1.00"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.01%  "

$ws.Range("D24").Value = "'0.516"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.95%  "

$ws.Range("E25").Value = "  +1.16%  "

$ws.Range("E26").Value = "  +3.99%  "

$ws.Range("D27").Value = "'9.65"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.52%  "

$ws.Range("D28").Value = "'0.997"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -1.03%  "

$ws.Range("E29").Value = "  +0.98%  "

$ws.Range("D30").Value = "'5.74"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.20%  "

$ws.Range("D31").Value = "'23.00"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.73%  "

$ws.Range("E32").Value = "  +2.42%  "

$ws.Range("E33").Value = "  +3.57%  "

$ws.Range("E34").Value = "  +0.05%  "

$ws.Range("B35").Value = "Monero"
$ws.Range("C35").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D35").Value = "'164.38"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.61%  "

$ws.Range("B36").Value = "ImmutableX"
$ws.Range("C36").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D36").Value = "'1.50"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.65%  "

$ws.Range("E37").Value = "  +1.32%  "

$ws.Range("E38").Value = "  -2.53%  "

$ws.Range("D39").Value = "'4.60"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +2.43%  "

$ws.Range("D40").Value = "'26.21"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -2.23%  "

$ws.Range("D41").Value = "'6.55"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -3.56%  "

$ws.Range("D42").Value = "'2.57"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.98%  "

$ws.Range("D43").Value = "'41.51"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +2.25%  "

$ws.Range("D44").Value = "'0.0688"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +2.05%  "

$ws.Range("D45").Value = "'342.31"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -5.29%  "

$ws.Range("D46").Value = "2.609.15"
$ws.Range("E46").Value = "  -5.09%  "

$ws.Range("D47").Value = "'24.81"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -2.47%  "

$ws.Range("D48").Value = "'0.0282"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.20%  "

$ws.Range("D49").Value = "'32.06"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +3.78%  "

$ws.Range("D50").Value = "'6.31"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +2.74%  "

$ws.Range("E51").Value = "  -0.20%  "
